# Auto-update stock values: 2025-12-12 07:55:10 UTC
# Adds one new trailing date column (20251211) to each of the 12 time-series
# sheets (시가/고가/저가/종가/거래량/s20/s60/z20/z60/gap/std/quant), copying the
# format of the previous last column and filling in the header date plus the
# two data rows.

$wb = $excel.ActiveWorkbook

function Add-DateColumn($SheetIndex, $NewCol, $IsText, $Row2, $Row3) {
    $ws = $wb.Worksheets.Item($SheetIndex)
    $prevCol = $NewCol - 1

    # Match the new column's width to the previous (last data) column, same
    # as every other column on the sheet.
    $ws.Columns.Item($NewCol).ColumnWidth = $ws.Columns.Item($prevCol).ColumnWidth

    # Copy the formatting (style, etc.) of the previous last column's header
    # cell onto the new header cell so it keeps the same bold/shaded header
    # style (cellXfs index) without minting a new style.
    $ws.Cells.Item(1, $prevCol).Copy()
    $ws.Cells.Item(1, $NewCol).PasteSpecial(-4122)

    if ($IsText) {
        # A few sheets (gap/std/quant) store the date header as literal text
        # rather than a number. Force text storage via NumberFormat, assign
        # the value, then re-paste the original header formatting so the
        # final cell style matches the neighboring header cells exactly.
        $ws.Cells.Item(1, $NewCol).NumberFormat = "@"
        $ws.Cells.Item(1, $NewCol).Value = "20251211"
        $ws.Cells.Item(1, $prevCol).Copy()
        $ws.Cells.Item(1, $NewCol).PasteSpecial(-4122)
    } else {
        $ws.Cells.Item(1, $NewCol).Value = 20251211
    }

    $ws.Cells.Item(2, $NewCol).Value = $Row2
    $ws.Cells.Item(3, $NewCol).Value = $Row3
}

# 시가 (Open) - sheet 2, new column BW (75)
Add-DateColumn 2 75 $false 623.8200000000001 55.63

# 고가 (High) - sheet 3, new column BW (75)
Add-DateColumn 3 75 $false 625.78 56.16

# 저가 (Low) - sheet 4, new column BW (75)
Add-DateColumn 4 75 $false 617.72 53.98

# 종가 (Close) - sheet 5, new column BW (75)
Add-DateColumn 5 75 $false 625.58 56.11

# 거래량 (Volume) - sheet 6, new column BW (75)
Add-DateColumn 6 75 $false 58272844 99683926

# s20 - sheet 7, new column BD (56)
Add-DateColumn 7 56 $false 95 17

# s60 - sheet 8, new column P (16)
Add-DateColumn 8 16 $false 80 13

# z20 - sheet 9, new column BD (56)
Add-DateColumn 9 56 $false 51 -22

# z60 - sheet 10, new column P (16)
Add-DateColumn 10 16 $false 60 -79

# gap - sheet 11, new column BD (56) - header stored as text
Add-DateColumn 11 56 $true 102 85

# std - sheet 12, new column AK (37) - header stored as text
Add-DateColumn 12 37 $true 6.19 6.94

# quant - sheet 13, new column P (16) - header stored as text
Add-DateColumn 13 16 $true 50 72

Write-Output "Added 20251211 column to all 12 data sheets"
